# NSMB - Most of last stage done!
# Extends the "V4" run table (sheet1 / tab "V4") with new rows of timing
# data (rows 203-216), adds two annotation notes in column G, and keeps
# the existing D-column "segment length" shared formula going for the
# new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row data -----------------------------------------------------
# Each entry: row, A-label (or $null), B (start frame), C (end frame)
$rows = @(
    @{ r = 203; a = "End Level";         b = 69906; c = 80786 },
    @{ r = 204; a = "Enter 8-Final";     b = 70786; c = 82543 },
    @{ r = 205; a = "1st Move";          b = 71017; c = 82799 },
    @{ r = 206; a = "Enter door";        b = 71467; c = 83249 },
    @{ r = 207; a = "Touch button";      b = 71725; c = 83509 },
    @{ r = 208; a = "Enter door";        b = 72016; c = 83800 },
    @{ r = 209; a = "Checkpoint 1657";   b = 72323; c = 84128 },
    @{ r = 210; a = "Enter door";        b = 72665; c = 84487 },
    @{ r = 211; a = "Enter door";        b = 73412; c = 85239 },
    @{ r = 212; a = "Enter door";        b = 73958; c = 85786 },
    @{ r = 213; a = $null;               b = 74154; c = 85983 },
    @{ r = 214; a = $null;               b = 74292; c = 86121 },
    @{ r = 215; a = $null;               b = 74361; c = 86190 },
    @{ r = 216; a = $null;               b = 74597; c = 86424 }
)

foreach ($row in $rows) {
    $r = $row.r
    if ($row.a -ne $null) {
        $ws.Range("A$r").Value = $row.a
    }
    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = $row.c
}

# Annotation notes in column G for the two approximate checkpoints
$ws.Range("G207").Value = "Note: not an absolute measure point"
$ws.Range("G209").Value = "Approx (camera angle diffs)"

# --- Extend the D-column shared formula down through row 216 ----------
# D90 originally carried the shared formula (si=7) covering D90:D202;
# re-applying it across the whole D90:D216 span keeps every row's
# "elapsed time" calculation consistent (IF(B>0,C-B,0)).
$ws.Range("D90:D216").FormulaR1C1 = "=IF(RC[-2] >  0,RC[-1]-RC[-2], 0)"

# --- Update the view: frozen header row, scrolled down to the new tail -
$ws.Range("B217").Select()
